# Corrects the "Fitness" column (column C) in Sheet1, rows 2-252,
# to reflect the recomputed best-fitness-so-far values produced by the
# corrected SA algorithm (run_6 log).
#
# Column layout: A = Run, B = Generation (0-based), C = Fitness
# For each row, the new Fitness value depends on the Generation (B) value:
#   Generation 0  - 13  -> 8045
#   Generation 14 - 25  -> 7701
#   Generation 26 - 250 -> 7651

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 252

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $generation = $ws.Cells.Item($row, 2).Value2

    if ($generation -le 13) {
        $newFitness = 8045
    } elseif ($generation -le 25) {
        $newFitness = 7701
    } else {
        $newFitness = 7651
    }

    $ws.Cells.Item($row, 3).Value = $newFitness
}

Write-Output "Updated fitness values for rows $firstRow to $lastRow"
